# Update 'want to go' count (column F) values per the source diff.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Cells.Item(2, 6).Value = 677
$ws.Cells.Item(3, 6).Value = 1498
$ws.Cells.Item(4, 6).Value = 3254
$ws.Cells.Item(6, 6).Value = 671
$ws.Cells.Item(7, 6).Value = 2232
$ws.Cells.Item(8, 6).Value = 484
$ws.Cells.Item(9, 6).Value = 409
$ws.Cells.Item(10, 6).Value = 236
$ws.Cells.Item(11, 6).Value = 129
$ws.Cells.Item(12, 6).Value = 318
$ws.Cells.Item(13, 6).Value = 1074
$ws.Cells.Item(14, 6).Value = 434
$ws.Cells.Item(15, 6).Value = 11
$ws.Cells.Item(16, 6).Value = 80
$ws.Cells.Item(17, 6).Value = 214
$ws.Cells.Item(18, 6).Value = 4494
$ws.Cells.Item(19, 6).Value = 3
$ws.Cells.Item(20, 6).Value = 1301
$ws.Cells.Item(21, 6).Value = 3406
$ws.Cells.Item(23, 6).Value = 95
$ws.Cells.Item(24, 6).Value = 175
$ws.Cells.Item(25, 6).Value = 3418
$ws.Cells.Item(26, 6).Value = 4965
$ws.Cells.Item(29, 6).Value = 545
$ws.Cells.Item(30, 6).Value = 3206
$ws.Cells.Item(33, 6).Value = 132
$ws.Cells.Item(34, 6).Value = 87
$ws.Cells.Item(35, 6).Value = 874
$ws.Cells.Item(36, 6).Value = 1162
$ws.Cells.Item(37, 6).Value = 1404
$ws.Cells.Item(39, 6).Value = 1332
$ws.Cells.Item(40, 6).Value = 848
$ws.Cells.Item(41, 6).Value = 13
$ws.Cells.Item(42, 6).Value = 806
$ws.Cells.Item(43, 6).Value = 494
$ws.Cells.Item(45, 6).Value = 294
$ws.Cells.Item(46, 6).Value = 62
$ws.Cells.Item(47, 6).Value = 150
$ws.Cells.Item(49, 6).Value = 3715

$ws = $wb.Worksheets.Item("演出")
$ws.Cells.Item(7, 6).Value = 999
$ws.Cells.Item(11, 6).Value = 8
$ws.Cells.Item(20, 6).Value = 49

$ws = $wb.Worksheets.Item("本地生活")
$ws.Cells.Item(2, 6).Value = 2138

$ws = $wb.Worksheets.Item("全部类型")
$ws.Cells.Item(2, 6).Value = 2138
$ws.Cells.Item(3, 6).Value = 677
$ws.Cells.Item(4, 6).Value = 1498
$ws.Cells.Item(5, 6).Value = 3254
$ws.Cells.Item(7, 6).Value = 671
$ws.Cells.Item(9, 6).Value = 2232
$ws.Cells.Item(10, 6).Value = 484
$ws.Cells.Item(11, 6).Value = 409
$ws.Cells.Item(12, 6).Value = 236
$ws.Cells.Item(13, 6).Value = 999
$ws.Cells.Item(14, 6).Value = 129
$ws.Cells.Item(15, 6).Value = 318
$ws.Cells.Item(16, 6).Value = 1074
$ws.Cells.Item(17, 6).Value = 434
$ws.Cells.Item(18, 6).Value = 11
$ws.Cells.Item(19, 6).Value = 214
$ws.Cells.Item(20, 6).Value = 4494
$ws.Cells.Item(21, 6).Value = 1301
$ws.Cells.Item(23, 6).Value = 3406
$ws.Cells.Item(24, 6).Value = 3422
$ws.Cells.Item(25, 6).Value = 4965
$ws.Cells.Item(28, 6).Value = 3206
$ws.Cells.Item(31, 6).Value = 132
$ws.Cells.Item(32, 6).Value = 87
$ws.Cells.Item(33, 6).Value = 875
$ws.Cells.Item(34, 6).Value = 1162
$ws.Cells.Item(35, 6).Value = 1404
$ws.Cells.Item(37, 6).Value = 1332
$ws.Cells.Item(39, 6).Value = 848
$ws.Cells.Item(40, 6).Value = 494
$ws.Cells.Item(41, 6).Value = 49
$ws.Cells.Item(44, 6).Value = 294
$ws.Cells.Item(46, 6).Value = 62
$ws.Cells.Item(47, 6).Value = 150
$ws.Cells.Item(49, 6).Value = 3715
